$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: D2 cleared, B2 and C2 populated
$ws.Range("D2").ClearContents()
$ws.Range("B2").Value = 28.980596743227572
$ws.Range("C2").Value = 26.166898802591774

# Row 3: B3 cleared, C3 updated
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 24.257865186880167

# Selection narrowed from B1:AY3 to B1:E3
$ws.Range("B1:E3").Select()
